$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 8.324917666666666
$ws.Range("H2").Value = 24.974753
$ws.Range("I2").Value = 0.8193616330571973
$ws.Range("J2").Value = 0.8193616330571972
$ws.Range("M2").Value = 31.61061466666667
$ws.Range("N2").Value = 94.831844
$ws.Range("O2").Value = 0.8860472269592234
$ws.Range("P2").Value = 0.8860472269592234
$ws.Range("Q2").Value = 263.1557644927258
$ws.Range("R2").Value = 2368.401880434532
$ws.Range("S2").Value = 0.7259931028471105
$ws.Range("T2").Value = 0.7259931028471104
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 8.324917666666666
$ws.Range("H3").Value = 24.974753
$ws.Range("I3").Value = 0.8193616330571973
$ws.Range("J3").Value = 0.8193616330571972
$ws.Range("O3").Value = 0.04688826274109129
$ws.Range("P3").Value = 0.04688826274109129
$ws.Range("Q3").Value = 13.92580017400756
$ws.Range("R3").Value = 125.332201566068
$ws.Range("S3").Value = 0.0384184435307555
$ws.Range("T3").Value = 0.03841844353075549
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 8.324917666666666
$ws.Range("H4").Value = 24.974753
$ws.Range("I4").Value = 0.8193616330571973
$ws.Range("J4").Value = 0.8193616330571972
$ws.Range("M4").Value = 2.392593
$ws.Range("N4").Value = 7.177778999999999
$ws.Range("O4").Value = 0.06706451029968528
$ws.Range("P4").Value = 0.06706451029968527
$ws.Range("Q4").Value = 19.918139734843
$ws.Range("R4").Value = 179.263257613587
$ws.Range("S4").Value = 0.05495008667933136
$ws.Range("T4").Value = 0.05495008667933134
$ws.Range("G5").Value = 1.835330666666667
$ws.Range("H5").Value = 5.505992
$ws.Range("I5").Value = 0.1806383669428028
$ws.Range("J5").Value = 0.1806383669428027
$ws.Range("M5").Value = 31.61061466666667
$ws.Range("N5").Value = 94.831844
$ws.Range("O5").Value = 0.8860472269592234
$ws.Range("P5").Value = 0.8860472269592234
$ws.Range("Q5").Value = 58.01593048991644
$ws.Range("R5").Value = 522.1433744092481
$ws.Range("S5").Value = 0.1600541241121131
$ws.Range("T5").Value = 0.160054124112113
$ws.Range("G6").Value = 1.835330666666667
$ws.Range("H6").Value = 5.505992
$ws.Range("I6").Value = 0.1806383669428028
$ws.Range("J6").Value = 0.1806383669428027
$ws.Range("O6").Value = 0.04688826274109129
$ws.Range("P6").Value = 0.04688826274109129
$ws.Range("Q6").Value = 3.070114221016889
$ws.Range("R6").Value = 27.631027989152
$ws.Range("S6").Value = 0.008469819210335796
$ws.Range("T6").Value = 0.008469819210335794
$ws.Range("G7").Value = 1.835330666666667
$ws.Range("H7").Value = 5.505992
$ws.Range("I7").Value = 0.1806383669428028
$ws.Range("J7").Value = 0.1806383669428027
$ws.Range("M7").Value = 2.392593
$ws.Range("N7").Value = 7.177778999999999
$ws.Range("O7").Value = 0.06706451029968528
$ws.Range("P7").Value = 0.06706451029968527
$ws.Range("Q7").Value = 4.391199305751999
$ws.Range("R7").Value = 39.520793751768
$ws.Range("S7").Value = 0.01211442362035393
$ws.Range("T7").Value = 0.01211442362035392
